# Update container / trough helper text on the "Groups" and "Individuals"
# sheets of the mactaquac-distribution template, per the commit
# "update container names in mactaquac parsers".

$wb = $excel.ActiveWorkbook

$wsGroups = $wb.Worksheets.Item("Groups")
$wsIndividuals = $wb.Worksheets.Item("Individuals")

# --- Individuals sheet: update the "Tank" / "Trough" helper/example text ---
$wsIndividuals.Range("M3").Value = "Optional. Containers fish were taken from.  E.g. LP1"
$wsIndividuals.Range("N3").Value = "Optional. Container fish were taken from.  E.g. TR1"

# --- Groups sheet: update the "Trough" helper/example text in Q3 ---
$wsGroups.Range("Q3").Value = "Enter trough names here if distributing from troughs. Eg.TR 4,TR 5,TR6"

# row 3 on Individuals grew taller to fit the new wrapped text
$wsIndividuals.Range("A3").EntireRow.RowHeight = 63.75

# --- Update which sheet/cell is active & selected ---
# Individuals was the active tab with A2 selected; now Groups is the active
# tab (with P3 selected) and Individuals keeps N4 selected when revisited.
$wsIndividuals.Range("N4").Select()
$wsGroups.Activate()
$wsGroups.Range("P3").Select()
